{"js": "// Replace each arithmetic-equation cell in the single 20x5 table with its\n// updated equation, matching the target revision. Cells are addressed\n// positionally (row, col) so the edit is independent of any text collisions\n// between old/new values elsewhere in the table. The existing run formatting\n// (font, size, paragraph alignment) is preserved because we replace text\n// inside the matched search range rather than rebuilding the paragraph.\n\nconst pairs = [\n  [0, 0, \"61+27=88\", \"43-22=21\"],\n  [0, 1, \"5+75=80\", \"69+18=87\"],\n  [0, 2, \"27-26=1\", \"42+56=98\"],\n  [0, 3, \"75-52=23\", \"23+22=45\"],\n  [0, 4, \"30+30=60\", \"20+31=51\"],\n  [1, 0, \"30+10=40\", \"71-20=51\"],\n  [1, 1, \"87-57=30\", \"97-25=72\"],\n  [1, 2, \"57-48=9\", \"9+59=68\"],\n  [1, 3, \"56+5=61\", \"8+72=80\"],\n  [1, 4, \"78-38=40\", \"64+12=76\"],\n  [2, 0, \"47-22=25\", \"92-44=48\"],\n  [2, 1, \"91-26=65\", \"81-19=62\"],\n  [2, 2, \"25+9=34\", \"1+42=43\"],\n  [2, 3, \"38+49=87\", \"68-26=42\"],\n  [2, 4, \"78+15=93\", \"0+80=80\"],\n  [3, 0, \"30+48=78\", \"80-15=65\"],\n  [3, 1, \"93+2=95\", \"9+40=49\"],\n  [3, 2, \"41-35=6\", \"75+23=98\"],\n  [3, 3, \"57-15=42\", \"11+62=73\"],\n  [3, 4, \"69+17=86\", \"81+6=87\"],\n  [4, 0, \"86-86=0\", \"21-20=1\"],\n  [4, 1, \"28+44=72\", \"1+49=50\"],\n  [4, 2, \"89-79=10\", \"80-69=11\"],\n  [4, 3, \"93-7=86\", \"6+0=6\"],\n  [4, 4, \"91-57=34\", \"68-39=29\"],\n  [5, 0, \"27+71=98\", \"13+59=72\"],\n  [5, 1, \"61+1=62\", \"30+6=36\"],\n  [5, 2, \"77+11=88\", \"56-41=15\"],\n  [5, 3, \"60-1=59\", \"4+85=89\"],\n  [5, 4, \"25+56=81\", \"47+27=74\"],\n  [6, 0, \"96-9=87\", \"52-20=32\"],\n  [6, 1, \"78-54=24\", \"35-5=30\"],\n  [6, 2, \"94-3=91\", \"17+35=52\"],\n  [6, 3, \"31+48=79\", \"50-46=4\"],\n  [6, 4, \"34+44=78\", \"63-6=57\"],\n  [7, 0, \"77-56=21\", \"34+16=50\"],\n  [7, 1, \"68+15=83\", \"29-12=17\"],\n  [7, 2, \"83-13=70\", \"93-62=31\"],\n  [7, 3, \"84-6=78\", \"51+8=59\"],\n  [7, 4, \"40+7=47\", \"7+0=7\"],\n  [8, 0, \"42+44=86\", \"81-2=79\"],\n  [8, 1, \"80-48=32\", \"87-42=45\"],\n  [8, 2, \"13-7=6\", \"43-28=15\"],\n  [8, 3, \"43+31=74\", \"95-6=89\"],\n  [8, 4, \"9+72=81\", \"4+23=27\"],\n  [9, 0, \"74-10=64\", \"36+18=54\"],\n  [9, 1, \"94-93=1\", \"89-85=4\"],\n  [9, 2, \"18-3=15\", \"92-35=57\"],\n  [9, 3, \"84-56=28\", \"69-37=32\"],\n  [9, 4, \"19+43=62\", \"75+3=78\"],\n  [10, 0, \"79-7=72\", \"4+10=14\"],\n  [10, 1, \"2+6=8\", \"78-56=22\"],\n  [10, 2, \"28-14=14\", \"74-32=42\"],\n  [10, 3, \"10+82=92\", \"96-55=41\"],\n  [10, 4, \"36+16=52\", \"86+10=96\"],\n  [11, 0, \"23+9=32\", \"37-20=17\"],\n  [11, 1, \"95-60=35\", \"24+32=56\"],\n  [11, 2, \"98-73=25\", \"96-15=81\"],\n  [11, 3, \"36+20=56\", \"67+17=84\"],\n  [11, 4, \"28+33=61\", \"83-79=4\"],\n  [12, 0, \"62-1=61\", \"0+21=21\"],\n  [12, 1, \"71-2=69\", \"25+73=98\"],\n  [12, 2, \"73-48=25\", \"52+21=73\"],\n  [12, 3, \"43-4=39\", \"58-32=26\"],\n  [12, 4, \"94-41=53\", \"81+16=97\"],\n  [13, 0, \"12-1=11\", \"71+25=96\"],\n  [13, 1, \"9+65=74\", \"49-36=13\"],\n  [13, 2, \"44+21=65\", \"26+0=26\"],\n  [13, 3, \"60-29=31\", \"60+6=66\"],\n  [13, 4, \"33+35=68\", \"79-44=35\"],\n  [14, 0, \"52+4=56\", \"48+43=91\"],\n  [14, 1, \"30+45=75\", \"85-72=13\"],\n  [14, 2, \"57+21=78\", \"60-40=20\"],\n  [14, 3, \"58+34=92\", \"26-12=14\"],\n  [14, 4, \"22+72=94\", \"0+85=85\"],\n  [15, 0, \"5+60=65\", \"15-12=3\"],\n  [15, 1, \"65-42=23\", \"88-56=32\"],\n  [15, 2, \"71+13=84\", \"47-1=46\"],\n  [15, 3, \"51-51=0\", \"17+78=95\"],\n  [15, 4, \"16+24=40\", \"92-14=78\"],\n  [16, 0, \"2+63=65\", \"87-1=86\"],\n  [16, 1, \"54-9=45\", \"68-58=10\"],\n  [16, 2, \"52+15=67\", \"85-12=73\"],\n  [16, 3, \"9+49=58\", \"99-60=39\"],\n  [16, 4, \"96-66=30\", \"83-29=54\"],\n  [17, 0, \"21+37=58\", \"54-1=53\"],\n  [17, 1, \"89-82=7\", \"95-0=95\"],\n  [17, 2, \"99-32=67\", \"93+3=96\"],\n  [17, 3, \"40-11=29\", \"81-43=38\"],\n  [17, 4, \"16+21=37\", \"33-27=6\"],\n  [18, 0, \"48+48=96\", \"54+1=55\"],\n  [18, 1, \"51-43=8\", \"79+18=97\"],\n  [18, 2, \"0+14=14\", \"33-12=21\"],\n  [18, 3, \"19+51=70\", \"16+53=69\"],\n  [18, 4, \"1+10=11\", \"79-43=36\"],\n  [19, 0, \"33+34=67\", \"72-69=3\"],\n  [19, 1, \"61-40=21\", \"0+11=11\"],\n  [19, 2, \"10-9=1\", \"45-5=40\"],\n  [19, 3, \"64-2=62\", \"75-73=2\"],\n  [19, 4, \"87-42=45\", \"20-1=19\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of pairs) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Expected text \"${oldText}\" not found in cell (${row}, ${col})`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-equation cell in the single 20x5 table with its\n# updated equation, matching the target revision. Cells are addressed\n# positionally (row, col) via Table.Cell(), so the edit does not depend on\n# any text collisions between old/new values elsewhere in the table.\n# Assigning Cell.Range.Text preserves the cell's existing run formatting\n# (font, size, paragraph alignment) because Word keeps the run's rPr when\n# only the text content of a fully-selected range is replaced.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$pairs = @(\n    @(1, 1, \"61+27=88\", \"43-22=21\"),\n    @(1, 2, \"5+75=80\", \"69+18=87\"),\n    @(1, 3, \"27-26=1\", \"42+56=98\"),\n    @(1, 4, \"75-52=23\", \"23+22=45\"),\n    @(1, 5, \"30+30=60\", \"20+31=51\"),\n    @(2, 1, \"30+10=40\", \"71-20=51\"),\n    @(2, 2, \"87-57=30\", \"97-25=72\"),\n    @(2, 3, \"57-48=9\", \"9+59=68\"),\n    @(2, 4, \"56+5=61\", \"8+72=80\"),\n    @(2, 5, \"78-38=40\", \"64+12=76\"),\n    @(3, 1, \"47-22=25\", \"92-44=48\"),\n    @(3, 2, \"91-26=65\", \"81-19=62\"),\n    @(3, 3, \"25+9=34\", \"1+42=43\"),\n    @(3, 4, \"38+49=87\", \"68-26=42\"),\n    @(3, 5, \"78+15=93\", \"0+80=80\"),\n    @(4, 1, \"30+48=78\", \"80-15=65\"),\n    @(4, 2, \"93+2=95\", \"9+40=49\"),\n    @(4, 3, \"41-35=6\", \"75+23=98\"),\n    @(4, 4, \"57-15=42\", \"11+62=73\"),\n    @(4, 5, \"69+17=86\", \"81+6=87\"),\n    @(5, 1, \"86-86=0\", \"21-20=1\"),\n    @(5, 2, \"28+44=72\", \"1+49=50\"),\n    @(5, 3, \"89-79=10\", \"80-69=11\"),\n    @(5, 4, \"93-7=86\", \"6+0=6\"),\n    @(5, 5, \"91-57=34\", \"68-39=29\"),\n    @(6, 1, \"27+71=98\", \"13+59=72\"),\n    @(6, 2, \"61+1=62\", \"30+6=36\"),\n    @(6, 3, \"77+11=88\", \"56-41=15\"),\n    @(6, 4, \"60-1=59\", \"4+85=89\"),\n    @(6, 5, \"25+56=81\", \"47+27=74\"),\n    @(7, 1, \"96-9=87\", \"52-20=32\"),\n    @(7, 2, \"78-54=24\", \"35-5=30\"),\n    @(7, 3, \"94-3=91\", \"17+35=52\"),\n    @(7, 4, \"31+48=79\", \"50-46=4\"),\n    @(7, 5, \"34+44=78\", \"63-6=57\"),\n    @(8, 1, \"77-56=21\", \"34+16=50\"),\n    @(8, 2, \"68+15=83\", \"29-12=17\"),\n    @(8, 3, \"83-13=70\", \"93-62=31\"),\n    @(8, 4, \"84-6=78\", \"51+8=59\"),\n    @(8, 5, \"40+7=47\", \"7+0=7\"),\n    @(9, 1, \"42+44=86\", \"81-2=79\"),\n    @(9, 2, \"80-48=32\", \"87-42=45\"),\n    @(9, 3, \"13-7=6\", \"43-28=15\"),\n    @(9, 4, \"43+31=74\", \"95-6=89\"),\n    @(9, 5, \"9+72=81\", \"4+23=27\"),\n    @(10, 1, \"74-10=64\", \"36+18=54\"),\n    @(10, 2, \"94-93=1\", \"89-85=4\"),\n    @(10, 3, \"18-3=15\", \"92-35=57\"),\n    @(10, 4, \"84-56=28\", \"69-37=32\"),\n    @(10, 5, \"19+43=62\", \"75+3=78\"),\n    @(11, 1, \"79-7=72\", \"4+10=14\"),\n    @(11, 2, \"2+6=8\", \"78-56=22\"),\n    @(11, 3, \"28-14=14\", \"74-32=42\"),\n    @(11, 4, \"10+82=92\", \"96-55=41\"),\n    @(11, 5, \"36+16=52\", \"86+10=96\"),\n    @(12, 1, \"23+9=32\", \"37-20=17\"),\n    @(12, 2, \"95-60=35\", \"24+32=56\"),\n    @(12, 3, \"98-73=25\", \"96-15=81\"),\n    @(12, 4, \"36+20=56\", \"67+17=84\"),\n    @(12, 5, \"28+33=61\", \"83-79=4\"),\n    @(13, 1, \"62-1=61\", \"0+21=21\"),\n    @(13, 2, \"71-2=69\", \"25+73=98\"),\n    @(13, 3, \"73-48=25\", \"52+21=73\"),\n    @(13, 4, \"43-4=39\", \"58-32=26\"),\n    @(13, 5, \"94-41=53\", \"81+16=97\"),\n    @(14, 1, \"12-1=11\", \"71+25=96\"),\n    @(14, 2, \"9+65=74\", \"49-36=13\"),\n    @(14, 3, \"44+21=65\", \"26+0=26\"),\n    @(14, 4, \"60-29=31\", \"60+6=66\"),\n    @(14, 5, \"33+35=68\", \"79-44=35\"),\n    @(15, 1, \"52+4=56\", \"48+43=91\"),\n    @(15, 2, \"30+45=75\", \"85-72=13\"),\n    @(15, 3, \"57+21=78\", \"60-40=20\"),\n    @(15, 4, \"58+34=92\", \"26-12=14\"),\n    @(15, 5, \"22+72=94\", \"0+85=85\"),\n    @(16, 1, \"5+60=65\", \"15-12=3\"),\n    @(16, 2, \"65-42=23\", \"88-56=32\"),\n    @(16, 3, \"71+13=84\", \"47-1=46\"),\n    @(16, 4, \"51-51=0\", \"17+78=95\"),\n    @(16, 5, \"16+24=40\", \"92-14=78\"),\n    @(17, 1, \"2+63=65\", \"87-1=86\"),\n    @(17, 2, \"54-9=45\", \"68-58=10\"),\n    @(17, 3, \"52+15=67\", \"85-12=73\"),\n    @(17, 4, \"9+49=58\", \"99-60=39\"),\n    @(17, 5, \"96-66=30\", \"83-29=54\"),\n    @(18, 1, \"21+37=58\", \"54-1=53\"),\n    @(18, 2, \"89-82=7\", \"95-0=95\"),\n    @(18, 3, \"99-32=67\", \"93+3=96\"),\n    @(18, 4, \"40-11=29\", \"81-43=38\"),\n    @(18, 5, \"16+21=37\", \"33-27=6\"),\n    @(19, 1, \"48+48=96\", \"54+1=55\"),\n    @(19, 2, \"51-43=8\", \"79+18=97\"),\n    @(19, 3, \"0+14=14\", \"33-12=21\"),\n    @(19, 4, \"19+51=70\", \"16+53=69\"),\n    @(19, 5, \"1+10=11\", \"79-43=36\"),\n    @(20, 1, \"33+34=67\", \"72-69=3\"),\n    @(20, 2, \"61-40=21\", \"0+11=11\"),\n    @(20, 3, \"10-9=1\", \"45-5=40\"),\n    @(20, 4, \"64-2=62\", \"75-73=2\"),\n    @(20, 5, \"87-42=45\", \"20-1=19\")\n)\n\nforeach ($p in $pairs) {\n    $row = $p[0]\n    $col = $p[1]\n    $oldText = $p[2]\n    $newText = $p[3]\n\n    $cell = $t.Cell($row, $col)\n    $range = $cell.Range\n    $current = $range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -ne $oldText) {\n        throw \"Cell ($row, $col): expected `\"$oldText`\" but found `\"$current`\"\"\n    }\n\n    $range.Text = $newText\n}\n"}
